# "some color palettes, end scene"
# Update the intro_2 closing line and append the new "end scene" rows
# (end_0, end_1, complete) to the Language/language.xlsx lookup sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# intro_2's value text changed ("he" -> "Robert", "move forward" -> "move on")
$ws.Range("B57").Value = "Excellent! With your guidance, Robert can finally move on."

# New end-scene rows
$ws.Range("A58").Value = "end_0"
$ws.Range("B58").Value = "After a long arduous journey, Robert is finally united with his family!"
$ws.Range("C58").Value = 5

$ws.Range("A59").Value = "end_1"
$ws.Range("B59").Value = "Thank you for playing!"
$ws.Range("C59").Value = 2

$ws.Range("A60").Value = "complete"
$ws.Range("B60").Value = "COMPLETE"

# Match the author's final selection/scroll state
$ws.Range("B60").Select()
